# Fix formatting; Add Wnioski
#
# 1. Update the cached "datetimeFigureOut" footer date (2023-12-15 -> 2023-12-31)
#    on the slide master and every slide layout.
# 2. Append a parenthetical remark to the "Cel pracy" slide's title.
# 3. Insert a new "Schemat blokowy" slide right before the existing
#    "Wnioski" slide (i.e. as the new slide 7, pushing "Wnioski" to slide 8).

$p = $ppt.ActivePresentation

# --- 1. Footer date fields -------------------------------------------------

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "2023-12-31"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2. "Cel pracy" title addition -----------------------------------------

$celPracySlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).Shapes.Item(1).TextFrame.TextRange.Text -eq "Cel pracy") {
        $celPracySlide = $p.Slides.Item($i)
        break
    }
}

$titleRange = $celPracySlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.InsertAfter(" (bardziej opisac dlaczego I po co)") | Out-Null

# --- 3. Insert "Schemat blokowy" slide before "Wnioski" ---------------------

$wnioskiIndex = 0
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).Shapes.Item(1).TextFrame.TextRange.Text -eq "Wnioski") {
        $wnioskiIndex = $i
        break
    }
}

$newSlide = $p.Slides.Add($wnioskiIndex, 2)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Schemat blokowy"

Write-Output "done"
